$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 649.8
$ws.Range("J19").Value = 542.2857
$ws.Range("L19").Value = 542.2857
$ws.Range("N19").Value = -892.2857
$ws.Range("H33").Value = 282
$ws.Range("I33").Value = 285.55554
$ws.Range("J33").Value = 250
$ws.Range("K33").Value = 285.55554
$ws.Range("L33").Value = 250
$ws.Range("M33").Value = -56.55554000000001
$ws.Range("N33").Value = -708
$ws.Range("H112").Value = 2318
$ws.Range("I112").Value = 1750
$ws.Range("J112").Value = 2367.3914
$ws.Range("K112").Value = 5250
$ws.Range("L112").Value = 7102.174199999999
$ws.Range("M112").Value = -4142
$ws.Range("N112").Value = -9318.174199999999
$ws.Range("H132").Value = 1383.6072
$ws.Range("I132").Value = 1293.3043
$ws.Range("K132").Value = 3879.9129
$ws.Range("M132").Value = -1349.9129
$ws.Range("H135").Value = 1270.2667
$ws.Range("I135").Value = 1111.0358
$ws.Range("K135").Value = 9999.322200000001
$ws.Range("M135").Value = -7464.322200000001
$ws.Range("H137").Value = 2923.3062
$ws.Range("I137").Value = 2084
$ws.Range("J137").Value = 3872.087
$ws.Range("K137").Value = 6252
$ws.Range("L137").Value = 11616.261
$ws.Range("M137").Value = -3702
$ws.Range("N137").Value = -16716.261
$ws.Range("H138").Value = 2367.41
$ws.Range("I138").Value = 1110.1111
$ws.Range("J138").Value = 2491.7583
$ws.Range("K138").Value = 3330.3333
$ws.Range("L138").Value = 7475.2749
$ws.Range("M138").Value = 1809.6667
$ws.Range("N138").Value = -17755.2749

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13673.761
$ws.Range("I32").Value = 8772.102000000001
$ws.Range("K32").Value = 8772.102000000001
$ws.Range("M32").Value = -8485.102000000001
$ws.Range("H61").Value = 4473.8335
$ws.Range("I61").Value = 2491.5334
$ws.Range("K61").Value = 2491.5334
$ws.Range("M61").Value = -2279.5334
$ws.Range("H88").Value = 4915814
$ws.Range("I88").Value = 12350
$ws.Range("J88").Value = 11920763
$ws.Range("K88").Value = 12350
$ws.Range("L88").Value = 11920763
$ws.Range("M88").Value = -11944
$ws.Range("N88").Value = -11921575
$ws.Range("H91").Value = 4915814
$ws.Range("I91").Value = 12350
$ws.Range("J91").Value = 11920763
$ws.Range("K91").Value = 12350
$ws.Range("L91").Value = 11920763
$ws.Range("M91").Value = -10946
$ws.Range("N91").Value = -11923571
$ws.Range("H112").Value = 40000
$ws.Range("J112").Value = 40000
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42954
$ws.Range("H132").Value = 4841.9536
$ws.Range("I132").Value = 4529.793
$ws.Range("K132").Value = 13589.379
$ws.Range("M132").Value = -11059.379
$ws.Range("H136").Value = 4473.8335
$ws.Range("I136").Value = 2491.5334
$ws.Range("K136").Value = 7474.600199999999
$ws.Range("M136").Value = -4924.600199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2148.0386
$ws.Range("I86").Value = 2384.7778
$ws.Range("J86").Value = 1615.375
$ws.Range("K86").Value = 2384.7778
$ws.Range("L86").Value = 1615.375
$ws.Range("M86").Value = -1261.7778
$ws.Range("N86").Value = -3861.375
$ws.Range("H89").Value = 2148.0386
$ws.Range("I89").Value = 2384.7778
$ws.Range("J89").Value = 1615.375
$ws.Range("K89").Value = 11923.889
$ws.Range("L89").Value = 8076.875
$ws.Range("M89").Value = -6307.888999999999
$ws.Range("N89").Value = -19308.875
$ws.Range("H94").Value = 6668538
$ws.Range("J94").Value = 28574682
$ws.Range("L94").Value = 28574682
$ws.Range("N94").Value = -28575584
$ws.Range("H126").Value = 50618.184
$ws.Range("J126").Value = 50618.184
$ws.Range("L126").Value = 50618.184
$ws.Range("H134").Value = 4463.1665
$ws.Range("I134").Value = 3470.389
$ws.Range("J134").Value = 7441.5
$ws.Range("K134").Value = 10411.167
$ws.Range("L134").Value = 22324.5
$ws.Range("M134").Value = -7876.167000000001
$ws.Range("N134").Value = -27394.5
$ws.Range("N126").Value = -60498.184

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49714.355
$ws.Range("J20").Value = 49714.355
$ws.Range("L20").Value = 49714.355
$ws.Range("N20").Value = -50186.355
$ws.Range("H22").Value = 859.6
$ws.Range("I22").Value = 733
$ws.Range("J22").Value = 1049.5
$ws.Range("K22").Value = 733
$ws.Range("L22").Value = 1049.5
$ws.Range("M22").Value = -383
$ws.Range("N22").Value = -1749.5
$ws.Range("H30").Value = 49714.355
$ws.Range("J30").Value = 49714.355
$ws.Range("L30").Value = 49714.355
$ws.Range("N30").Value = -49896.355
$ws.Range("H31").Value = 4663.607
$ws.Range("I31").Value = 3974.8
$ws.Range("J31").Value = 5046.278
$ws.Range("K31").Value = 3974.8
$ws.Range("L31").Value = 5046.278
$ws.Range("M31").Value = -3679.8
$ws.Range("N31").Value = -5636.278
$ws.Range("H34").Value = 4663.607
$ws.Range("I34").Value = 3974.8
$ws.Range("J34").Value = 5046.278
$ws.Range("K34").Value = 3974.8
$ws.Range("L34").Value = 5046.278
$ws.Range("M34").Value = -3772.8
$ws.Range("N34").Value = -5450.278
$ws.Range("H128").Value = 49714.355
$ws.Range("J128").Value = 49714.355
$ws.Range("L128").Value = 49714.355
$ws.Range("N128").Value = -59674.355
$ws.Range("H132").Value = 4132.278
$ws.Range("I132").Value = 3492.4666
$ws.Range("J132").Value = 7331.3335
$ws.Range("K132").Value = 10477.3998
$ws.Range("L132").Value = 21994.0005
$ws.Range("M132").Value = -7947.399800000001
$ws.Range("N132").Value = -27054.0005
$ws.Range("H140").Value = 71749.914
$ws.Range("J140").Value = 72818.17999999999
$ws.Range("L140").Value = 72818.17999999999
$ws.Range("N140").Value = -83178.17999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1195.0646
$ws.Range("I122").Value = 1084.5
$ws.Range("J122").Value = 1233.5217
$ws.Range("K122").Value = 9760.5
$ws.Range("L122").Value = 11101.6953
$ws.Range("M122").Value = -7310.5
$ws.Range("N122").Value = -16001.6953
$ws.Range("H129").Value = 1623.2759
$ws.Range("J129").Value = 1810.7273
$ws.Range("L129").Value = 5432.1819
$ws.Range("N129").Value = -15432.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("H70").Value = 99748.836
$ws.Range("I70").Value = 164930.14
$ws.Range("J70").Value = 8495
$ws.Range("K70").Value = 164930.14
$ws.Range("L70").Value = 8495
$ws.Range("M70").Value = -164660.14
$ws.Range("N70").Value = -9035
$ws.Range("H73").Value = 99748.836
$ws.Range("I73").Value = 164930.14
$ws.Range("J73").Value = 8495
$ws.Range("K73").Value = 164930.14
$ws.Range("L73").Value = 8495
$ws.Range("M73").Value = -163994.14
$ws.Range("N73").Value = -10367
$ws.Range("H118").Value = 35714.145
$ws.Range("J118").Value = 35714.145
$ws.Range("L118").Value = 35714.145
$ws.Range("N118").Value = -39028.145
$ws.Range("H132").Value = 4292.347
$ws.Range("I132").Value = 3544.5122
$ws.Range("J132").Value = 8125
$ws.Range("K132").Value = 10633.5366
$ws.Range("L132").Value = 24375
$ws.Range("M132").Value = -8103.536599999999
$ws.Range("N132").Value = -29435
$ws.Range("N27").ClearContents()  # removes cell entirely (was -5082)

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 425
$ws.Range("I16").Value = 425
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 425
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = -255
$ws.Range("H22").Value = 1252.6666
$ws.Range("I22").Value = 1263.4546
$ws.Range("K22").Value = 1263.4546
$ws.Range("M22").Value = -968.4546
$ws.Range("H27").Value = 1252.6666
$ws.Range("I27").Value = 1263.4546
$ws.Range("K27").Value = 1263.4546
$ws.Range("M27").Value = -1156.4546
$ws.Range("H61").Value = 2588.9443
$ws.Range("I61").Value = 1686.1786
$ws.Range("J61").Value = 5748.625
$ws.Range("K61").Value = 1686.1786
$ws.Range("L61").Value = 5748.625
$ws.Range("M61").Value = -1484.1786
$ws.Range("N61").Value = -6152.625
$ws.Range("H68").Value = 5810.231
$ws.Range("I68").Value = 5251.1763
$ws.Range("J68").Value = 6866.222
$ws.Range("K68").Value = 5251.1763
$ws.Range("L68").Value = 6866.222
$ws.Range("M68").Value = -4502.1763
$ws.Range("N68").Value = -8364.222
$ws.Range("H71").Value = 5810.231
$ws.Range("I71").Value = 5251.1763
$ws.Range("J71").Value = 6866.222
$ws.Range("K71").Value = 26255.8815
$ws.Range("L71").Value = 34331.11
$ws.Range("M71").Value = -22511.8815
$ws.Range("N71").Value = -41819.11
$ws.Range("H113").Value = 2588.9443
$ws.Range("I113").Value = 1686.1786
$ws.Range("J113").Value = 5748.625
$ws.Range("K113").Value = 1686.1786
$ws.Range("L113").Value = 5748.625
$ws.Range("M113").Value = 483.8214
$ws.Range("N113").Value = -10088.625
$ws.Range("H127").Value = 70000
$ws.Range("J127").Value = 70000
$ws.Range("L127").Value = 70000
$ws.Range("N127").Value = -79920
$ws.Range("H132").Value = 4960.5
$ws.Range("I132").Value = 4282.4287
$ws.Range("J132").Value = 6255
$ws.Range("K132").Value = 12847.2861
$ws.Range("L132").Value = 18765
$ws.Range("M132").Value = -10317.2861
$ws.Range("N132").Value = -23825
$ws.Range("H136").Value = 4680.931
$ws.Range("I136").Value = 3279.85
$ws.Range("J136").Value = 7794.4443
$ws.Range("K136").Value = 9839.549999999999
$ws.Range("L136").Value = 23383.3329
$ws.Range("M136").Value = -7289.549999999999
$ws.Range("N136").Value = -28483.3329
$ws.Range("N16").ClearContents()  # removes cell entirely (was -2315)

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 102569
$ws.Range("I96").Value = 127339.25
$ws.Range("J96").Value = 3488
$ws.Range("K96").Value = 127339.25
$ws.Range("L96").Value = 3488
$ws.Range("M96").Value = -125966.25
$ws.Range("N96").Value = -6234
$ws.Range("H132").Value = 1930.0227
$ws.Range("I132").Value = 1187.0541
$ws.Range("K132").Value = 3561.1623
$ws.Range("M132").Value = -1031.1623
